$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: UV Aerosol Index dataset
$ws.Hyperlinks.Add($ws.Range("C7"), "https://disc.gsfc.nasa.gov/datasets/OMPS_NPP_NMTO3_L3_DAILY_2/summary", "", "", "https://disc.gsfc.nasa.gov/datasets/OMPS_NPP_NMTO3_L3_DAILY_2/summary")
$ws.Range("A7").Value = "UV Aerosol Index - OMPS-NPP NMTO3 L2 v2.1"
$ws.Range("B7").Value = "NASA GES DIC"
$ws.Range("C7").Style = $ws.Range("C4").Style

# Row 8: Hourly Cloud Optical Depth dataset
$ws.Range("A8").Value = "Hourly Cloud Optical Depth - SATCORPS CERES GEO GOES-15 Ed 4 v1.0, v1.2 "
$ws.Range("B8").Value = "NOAA/NASA LARC DAP"
$ws.Hyperlinks.Add($ws.Range("C8"), "https://opendap.larc.nasa.gov/opendap/", "", "", "https://opendap.larc.nasa.gov/opendap/")
$ws.Range("C8").Style = $ws.Range("C4").Style

$ws.Range("D8").Select()
